$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the two rows that held the placeholder/test tracking numbers
# ("2435435" was row 2, "34534543" was row 4). Deleting row 2 first shifts
# the former row 4 up to row 3, so delete that row next.
$ws.Rows("2:2").Delete()
$ws.Rows("3:3").Delete()

# Renumber the "序号" (serial number) column sequentially for the
# remaining 13 tracking-number rows.
for ($i = 0; $i -lt 13; $i++) {
    $ws.Cells.Item(2 + $i, 1).Value = $i + 1
}

$ws.Range("E15").Select()
